# Edit LOB1214.docx:
#  1. In the "Avaliacao" section, the "Metodo" description run is split into two
#     lines (a manual line break is inserted after "...exercicios dirigidos. ").
#  2. In the "Bibliografia" section, the single run containing every reference
#     back-to-back is split so that each reference sits on its own line (manual
#     line breaks between references, and three line breaks before the
#     "Bibliografia complementar:" heading).

$d = $word.ActiveDocument

# 1. Split the "Metodo" run's text into two runs separated by a manual line break.
$found1 = $d.Content.Find.Execute(
    "Aulas teóricas e práticas, visitas técnicas e exercícios dirigidos. Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Aulas teóricas e práticas, visitas técnicas e exercícios dirigidos. ^lAvaliação baseada em provas, exercícios e trabalhos práticos e relatórios.",
    2)
Write-Host "Metodo split applied:" $found1
if (-not $found1) { throw "Could not find the 'Metodo' paragraph text to split." }

# 2. Split the bibliography run's text into individual runs separated by manual line breaks
#    (one per reference), with three line breaks before "Bibliografia complementar:".
$found2 = $d.Content.Find.Execute(
    "BURROUGH, P. A. Principles of Geographical Information Systems - Spatial Information Systems and Geoestatistics, Oxford: Clarendon Press, 1998.BURROUGH, P. A.; MCDONNELL, R. A. Principles of Geographical Information Systems. Oxford University Press, 1998.CÂMARA, G. & MEDEIROS, J. S. GIS para Meio Ambiente. INPE. São José dos Campos, SP, 1998.CROSTA, A. P. Processamento Digital de Imagens de Sensoriamento Remoto. Campinas – SP, 1992.FLORENZANO, T. G. Imagens de Satélite para Estudos Ambientais. Oficina de textos. São Paulo, 2002.IBGE. Noções Básicas de Cartografia. Rio de Janeiro. Coleção Manuais Técnicos em Geociências, 1999.LONGLEY, P. A.; GOODCHILD, M. F.; MAGUIRE, D. J.; RHIND, D. W. Geographic Information Systems and Science. John Wiley & Sons, 2001.MIRANDA, J. I.; Fundamentos de Sistemas de Informações Geográficas. Brasília, Embrapa, 2005.MOREIRA, M. A. Fundamentos do Sensoriamento Remoto e Metodologias de Aplicação. São José dos Campos – SP – INPE, 2001.SILVA, A.B. Sistemas de Informações Geo-referenciadas. Editora da Unicamp. Campinas. 1999.SILVA, A. B; Sistemas de informações Geo-referenciadas: conceitos e fundamentos. Campinas: Editora da Unicamp, 2003.SILVA, J.X. Geoprocessamento para Análise Ambiental. Rio de Janeiro. 2001.Bibliografia complementar:CARVALHO, M. S.; PINA, M. F.; SANTOS, S. M.  Conceitos Básicos de Sistemas de Informação Geográfica e Cartografia Aplicados à Saúde. Rede Interagencial de Informações para a Saúde. Brasília. Ministério da Saúde, 2000.DENT, B. D.  Cartography Thematic Map Design. 5th Edition. WCB/McGraw-Hill, 1999.MATOS, J. Fundamentos da Informação Geográfica. Lisboa, Lidel, 2008.MORAES NOVO, E. M. L. Sensoriamento Remoto – Princípios e Aplicações. 2ªEdição. São Paulo, 1992.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "BURROUGH, P. A. Principles of Geographical Information Systems - Spatial Information Systems and Geoestatistics, Oxford: Clarendon Press, 1998.^lBURROUGH, P. A.; MCDONNELL, R. A. Principles of Geographical Information Systems. Oxford University Press, 1998.^lCÂMARA, G. & MEDEIROS, J. S. GIS para Meio Ambiente. INPE. São José dos Campos, SP, 1998.^lCROSTA, A. P. Processamento Digital de Imagens de Sensoriamento Remoto. Campinas – SP, 1992.^lFLORENZANO, T. G. Imagens de Satélite para Estudos Ambientais. Oficina de textos. São Paulo, 2002.^lIBGE. Noções Básicas de Cartografia. Rio de Janeiro. Coleção Manuais Técnicos em Geociências, 1999.^lLONGLEY, P. A.; GOODCHILD, M. F.; MAGUIRE, D. J.; RHIND, D. W. Geographic Information Systems and Science. John Wiley & Sons, 2001.^lMIRANDA, J. I.; Fundamentos de Sistemas de Informações Geográficas. Brasília, Embrapa, 2005.^lMOREIRA, M. A. Fundamentos do Sensoriamento Remoto e Metodologias de Aplicação. São José dos Campos – SP – INPE, 2001.^lSILVA, A.B. Sistemas de Informações Geo-referenciadas. Editora da Unicamp. Campinas. 1999.^lSILVA, A. B; Sistemas de informações Geo-referenciadas: conceitos e fundamentos. Campinas: Editora da Unicamp, 2003.^lSILVA, J.X. Geoprocessamento para Análise Ambiental. Rio de Janeiro. 2001.^l^l^lBibliografia complementar:^lCARVALHO, M. S.; PINA, M. F.; SANTOS, S. M.  Conceitos Básicos de Sistemas de Informação Geográfica e Cartografia Aplicados à Saúde. Rede Interagencial de Informações para a Saúde. Brasília. Ministério da Saúde, 2000.^lDENT, B. D.  Cartography Thematic Map Design. 5th Edition. WCB/McGraw-Hill, 1999.^lMATOS, J. Fundamentos da Informação Geográfica. Lisboa, Lidel, 2008.^lMORAES NOVO, E. M. L. Sensoriamento Remoto – Princípios e Aplicações. 2ªEdição. São Paulo, 1992.",
    2)
Write-Host "Bibliografia split applied:" $found2
if (-not $found2) { throw "Could not find the 'Bibliografia' paragraph text to split." }
